$d = $word.ActiveDocument

# --------------------------------------------------------------------------
# Edit 1: Insert six spaces after "Christopher Landry" (same bold / sz32 run
#         formatting) and place the "_GoBack" bookmark right after them.
# --------------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.ClearFormatting()
$rng1.Find.Execute("Christopher Landry") | Out-Null
$rng1.Collapse(0)
$rng1.InsertAfter("      ")
$rng1.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rng1)

# --------------------------------------------------------------------------
# Edit 2: Merge the split hyperlink text "christ" + "o" + "p" + "h" + "erlandry.ca"
#         into a single run reading "christopherlandry.ca" (keeping the
#         Hyperlink character style on the resulting run).
# --------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$rng2.Find.Execute("christ" + "o" + "p" + "h" + "erlandry.ca") | Out-Null
$start2 = $rng2.Start
$rng2.Text = "@@SENTINEL2@@"
$rng2b = $d.Range($start2, $start2 + 13)
$rng2b.Text = "christopherlandry.ca"
$rng2b.Style = "Hyperlink"

# --------------------------------------------------------------------------
# Edit 3: Merge "Provided support for inbound and outbou" + the old "_GoBack"
#         bookmark + "nd billing " into one run, removing that old bookmark.
# --------------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.ClearFormatting()
$rng3.Find.Execute("Provided support for inbound and outbou" + "nd billing ") | Out-Null
$start3 = $rng3.Start
$rng3.Text = "@@SENTINEL3@@"
$rng3b = $d.Range($start3, $start3 + 13)
$rng3b.Text = "Provided support for inbound and outbound billing "

# --------------------------------------------------------------------------
# Edit 4: "Saint Lawrence College Research & Development" -> "St. Lawrence College R&D"
# --------------------------------------------------------------------------
$rng4 = $d.Content
$rng4.Find.ClearFormatting()
$rng4.Find.Execute("Saint Lawrence College Research & Development") | Out-Null
$start4 = $rng4.Start
$rng4.Text = "@@SENTINEL4@@"
$rng4b = $d.Range($start4, $start4 + 13)
$rng4b.Text = "St. Lawrence College R&D"
